$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Edit 1: shape "五边形 108" (Id=109) -----------------------------------
# "...中出现过的记录" -> "...中出现过的记录（中途转分队导致）"
# split into 5 runs (identical formatting) as in the target OOXML.
$sh1 = $s.Shapes.Item(5)
$tr1 = $sh1.TextFrame.TextRange

$target = $tr1.Characters(14, 7)   # "中出现过的记录"
$target.Text = "中出现过的记录（中途转分队导致）"

$r2 = $tr1.Characters(19, 2)       # "记录"
$r2.Font.Size = 9
$r3 = $tr1.Characters(21, 6)       # "（中途转分队"
$r3.Font.Size = 9
$r4 = $tr1.Characters(27, 2)       # "导致"
$r4.Font.Size = 9
$r5 = $tr1.Characters(29, 1)       # "）"
$r5.Font.Size = 9

# --- Edit 2: shape "圆角矩形 101" (Id=102) ----------------------------------
# "本周" + "所有成员" -> single run "本周所有成员"
$sh2 = $s.Shapes.Item(51)
$tr2 = $sh2.TextFrame.TextRange
$merge = $tr2.Characters(4, 6)     # "本周所有成员"
$merge.Text = "本周所有成员"
